# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos sheet
# with the latest scraped values. Price cells that look like plain numbers
# are written with a leading apostrophe so Excel keeps them as text (matching
# the workbook's existing text-formatted Price column) instead of coercing
# them into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.042.13'
$ws.Range("E2").Value = '  +2.92%  '
$ws.Range("D3").Value = '1.597.01'
$ws.Range("E3").Value = '  +1.98%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''212.42'
$ws.Range("E5").Value = '  +2.48%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").Value = '''0.484'
$ws.Range("E7").Value = '  +1.45%  '
$ws.Range("E8").Value = '  +2.48%  '
$ws.Range("E9").Value = '  +1.36%  '
$ws.Range("D10").Value = '''17.93'
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("D11").Value = '''0.0815'
$ws.Range("E11").Value = '  +4.42%  '
$ws.Range("D12").Value = '1.819.57'
$ws.Range("E12").Value = '  +1.98%  '
$ws.Range("D13").Value = '1.594.18'
$ws.Range("E13").Value = '  +1.70%  '
$ws.Range("E14").Value = '  -0.56%  '
$ws.Range("D15").Value = '''0.511'
$ws.Range("E15").Value = '  +0.86%  '
$ws.Range("D16").Value = '26.012.78'
$ws.Range("E16").Value = '  +2.77%  '
$ws.Range("D17").Value = '''60.33'
$ws.Range("E17").Value = '  +1.51%  '
$ws.Range("E18").Value = '  +1.21%  '
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("D20").Value = '''203.36'
$ws.Range("E20").Value = '  +9.70%  '
$ws.Range("D21").Value = '''4.23'
$ws.Range("E21").Value = '  +2.44%  '
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").Value = '''5.97'
$ws.Range("E23").Value = '  +1.50%  '
$ws.Range("E24").Value = '  +11.65%  '
$ws.Range("D25").Value = '''141.23'
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("E27").Value = '  -2.40%  '
$ws.Range("E28").Value = '  +2.50%  '
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("E30").Value = '  +1.17%  '
$ws.Range("E31").Value = '  +1.32%  '
$ws.Range("D32").Value = '''3.12'
$ws.Range("E32").Value = '  +2.95%  '
$ws.Range("E33").Value = '  -1.55%  '
$ws.Range("E34").Value = '  +0.85%  '
$ws.Range("E35").Value = '  +1.90%  '
$ws.Range("D36").Value = '1.108.78'
$ws.Range("E37").Value = '  +7.98%  '
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("D39").Value = '''2.32'
$ws.Range("E39").Value = '  +0.57%  '
$ws.Range("E40").Value = '  +0.97%  '
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("E42").Value = '  -4.43%  '
$ws.Range("D43").Value = '1.732.18'
$ws.Range("E43").Value = '  +1.97%  '
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("D45").Value = '''92.43'
$ws.Range("E45").Value = '  -0.79%  '
$ws.Range("D46").Value = '''1.49'
$ws.Range("E46").Value = '  +3.86%  '
$ws.Range("D47").Value = '''53.35'
$ws.Range("E47").Value = '  +1.44%  '
$ws.Range("E48").Value = '  +0.04%  '
$ws.Range("E49").Value = '  +0.64%  '
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("D51").Value = '0.0₇0926'
$ws.Range("E51").Value = '  -17.33%  '
